# Updates the Price (D) and Volume(1h) (E) columns of the crypto tracker sheet
# with freshly scraped values, written back as plain text (matching the
# original cell formatting, which stores these as text, not numbers/percentages).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "320.11" },
    @{ Cell = "E2"; Value = "3.82%" },
    @{ Cell = "D3"; Value = "41.33" },
    @{ Cell = "E3"; Value = "1.29%" },
    @{ Cell = "D4"; Value = "5.237" },
    @{ Cell = "E4"; Value = "2.16%" },
    @{ Cell = "D5"; Value = "0.07729" },
    @{ Cell = "E5"; Value = "1.49%" },
    @{ Cell = "D6"; Value = "1.699" },
    @{ Cell = "E6"; Value = "5.17%" },
    @{ Cell = "D7"; Value = "0.9459" },
    @{ Cell = "E7"; Value = "4.22%" },
    @{ Cell = "E8"; Value = "-1.30%" },
    @{ Cell = "D9"; Value = "0.1243" },
    @{ Cell = "E9"; Value = "-2.55%" },
    @{ Cell = "D10"; Value = "0.1827" },
    @{ Cell = "E10"; Value = "1.04%" },
    @{ Cell = "D11"; Value = "0.09175" },
    @{ Cell = "E11"; Value = "1.11%" },
    @{ Cell = "D12"; Value = "0.04339" },
    @{ Cell = "E12"; Value = "0.61%" },
    @{ Cell = "D13"; Value = "0.1052" },
    @{ Cell = "E13"; Value = "0.72%" },
    @{ Cell = "D14"; Value = "0.001293" },
    @{ Cell = "E14"; Value = "3.19%" },
    @{ Cell = "D15"; Value = "0.006023" },
    @{ Cell = "E15"; Value = "3.33%" },
    @{ Cell = "D17"; Value = "3.341" },
    @{ Cell = "E17"; Value = "-0.25%" },
    @{ Cell = "D18"; Value = "4.342" },
    @{ Cell = "E18"; Value = "1.47%" },
    @{ Cell = "D19"; Value = "0.3357" },
    @{ Cell = "E19"; Value = "1.31%" },
    @{ Cell = "D20"; Value = "7.740" },
    @{ Cell = "E20"; Value = "11.28%" },
    @{ Cell = "D21"; Value = "0.1354" },
    @{ Cell = "E21"; Value = "-2.81%" },
    @{ Cell = "D22"; Value = "0.2826" },
    @{ Cell = "E22"; Value = "4.37%" },
    @{ Cell = "D23"; Value = "0.04042" },
    @{ Cell = "E23"; Value = "-0.10%" },
    @{ Cell = "D24"; Value = "0.001266" },
    @{ Cell = "D25"; Value = "0.004111" },
    @{ Cell = "E25"; Value = "1.93%" },
    @{ Cell = "D26"; Value = "0.0001270" },
    @{ Cell = "E26"; Value = "-0.25%" },
    @{ Cell = "D38"; Value = "0.02544" },
    @{ Cell = "E38"; Value = "4.68%" },
    @{ Cell = "D39"; Value = "0.05337" },
    @{ Cell = "E39"; Value = "1.82%" },
    @{ Cell = "D40"; Value = "0.007773" },
    @{ Cell = "E40"; Value = "-0.94%" },
    @{ Cell = "D41"; Value = "0.1316" },
    @{ Cell = "E41"; Value = "1.51%" },
    @{ Cell = "D42"; Value = "0.007363" },
    @{ Cell = "E42"; Value = "8.14%" },
    @{ Cell = "D43"; Value = "0.001993" },
    @{ Cell = "E43"; Value = "4.88%" },
    @{ Cell = "D44"; Value = "0.008367" },
    @{ Cell = "E44"; Value = "13.43%" },
    @{ Cell = "D45"; Value = "0.3180" },
    @{ Cell = "E45"; Value = "-4.84%" },
    @{ Cell = "D46"; Value = "0.00006708" },
    @{ Cell = "E46"; Value = "-2.84%" },
    @{ Cell = "E47"; Value = "-0.24%" },
    @{ Cell = "D48"; Value = "0.2025" },
    @{ Cell = "E48"; Value = "89.88%" },
    @{ Cell = "D49"; Value = "0.004206" },
    @{ Cell = "E49"; Value = "40.08%" },
    @{ Cell = "D50"; Value = "0.00002100" },
    @{ Cell = "E50"; Value = "-0.24%" },
    @{ Cell = "D51"; Value = "0.0002000" },
    @{ Cell = "E51"; Value = "-0.24%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so Excel does not reinterpret "320.11" or "3.82%"
    # as a number/percentage (the source cells are plain text strings).
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Drop the now-unneeded text format so the cell style matches the
    # original (unstyled) cell.
    $cell.ClearFormats()
}
